$wb = $excel.ActiveWorkbook

# Sheet ALC, row 129 (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1124.921
$ws.Range("I129").Value = 799.75
$ws.Range("J129").Value = 1163.1765
$ws.Range("K129").Value = 2399.25
$ws.Range("L129").Value = 3489.5295
$ws.Range("M129").Value = 2600.75
$ws.Range("N129").Value = -13489.5295

# Sheet ARM, row 62 (hunk 1)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 22499.5
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 34999
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 34999
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -36247

# Sheet ARM, row 65 (hunk 2)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 22499.5
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 34999
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 104997
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -111237

# Sheet ARM, row 76 (hunk 3)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 13939.429
$ws.Range("J76").Value = 13939.429
$ws.Range("L76").Value = 13939.429
$ws.Range("N76").Value = -14615.429

# Sheet ARM, row 79 (hunk 4)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 13939.429
$ws.Range("J79").Value = 13939.429
$ws.Range("L79").Value = 13939.429
$ws.Range("N79").Value = -16279.429

# Sheet ARM, row 92 (hunk 5)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 19765.445
$ws.Range("J92").Value = 19765.445
$ws.Range("L92").Value = 19765.445
$ws.Range("N92").Value = -24757.445

# Sheet ARM, row 97 (hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1358.7742
$ws.Range("I97").Value = 1219.5834
$ws.Range("K97").Value = 1219.5834
$ws.Range("M97").Value = -723.5834

# Sheet BSM, row 76 (hunk 7)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 19500
$ws.Range("J76").Value = 19500
$ws.Range("L76").Value = 19500
$ws.Range("N76").Value = -20130

# Sheet BSM, row 79 (hunk 8)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 19500
$ws.Range("J79").Value = 19500
$ws.Range("L79").Value = 19500
$ws.Range("N79").Value = -21684

# Sheet BSM, row 86 (hunk 9)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2433.2666
$ws.Range("I86").Value = 2679.9
$ws.Range("J86").Value = 1940
$ws.Range("K86").Value = 2679.9
$ws.Range("L86").Value = 1940
$ws.Range("M86").Value = -1556.9
$ws.Range("N86").Value = -4186

# Sheet BSM, row 89 (hunk 10)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2433.2666
$ws.Range("I89").Value = 2679.9
$ws.Range("J89").Value = 1940
$ws.Range("K89").Value = 13399.5
$ws.Range("L89").Value = 9700
$ws.Range("M89").Value = -7783.5
$ws.Range("N89").Value = -20932

# Sheet CRP, row 31 (hunk 11)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1696.7551
$ws.Range("I31").Value = 1311.3077
$ws.Range("J31").Value = 3200
$ws.Range("K31").Value = 1311.3077
$ws.Range("L31").Value = 3200
$ws.Range("M31").Value = -1016.3077
$ws.Range("N31").Value = -3790

# Sheet CRP, row 34 (hunk 12)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1696.7551
$ws.Range("I34").Value = 1311.3077
$ws.Range("J34").Value = 3200
$ws.Range("K34").Value = 1311.3077
$ws.Range("L34").Value = 3200
$ws.Range("M34").Value = -1109.3077
$ws.Range("N34").Value = -3604

# Sheet CRP, row 58 (hunk 13)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1510.2972
$ws.Range("I58").Value = 746.2632
$ws.Range("J58").Value = 2316.7778
$ws.Range("K58").Value = 746.2632
$ws.Range("L58").Value = 2316.7778
$ws.Range("M58").Value = -543.2632
$ws.Range("N58").Value = -2722.7778

# Sheet CRP, row 62 (hunk 14)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4367.615
$ws.Range("I62").Value = 2534.8
$ws.Range("J62").Value = 5513.125
$ws.Range("K62").Value = 2534.8
$ws.Range("L62").Value = 5513.125
$ws.Range("M62").Value = -1910.8
$ws.Range("N62").Value = -6761.125

# Sheet CRP, row 65 (hunk 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4367.615
$ws.Range("I65").Value = 2534.8
$ws.Range("J65").Value = 5513.125
$ws.Range("K65").Value = 12674
$ws.Range("L65").Value = 27565.625
$ws.Range("M65").Value = -9554
$ws.Range("N65").Value = -33805.625

# Sheet CRP, row 132 (hunk 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2559.44
$ws.Range("I132").Value = 1798.7
$ws.Range("K132").Value = 5396.1
$ws.Range("M132").Value = -2866.1

# Sheet CRP, row 136 (hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1510.2972
$ws.Range("I136").Value = 746.2632
$ws.Range("J136").Value = 2316.7778
$ws.Range("K136").Value = 2238.7896
$ws.Range("L136").Value = 6950.3334
$ws.Range("M136").Value = 311.2103999999999
$ws.Range("N136").Value = -12050.3334

# Sheet CUL, row 12 (hunk 18)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 113.09091
$ws.Range("I12").Value = 3.75
$ws.Range("J12").Value = 175.57143
$ws.Range("K12").Value = 11.25
$ws.Range("L12").Value = 526.71429
$ws.Range("M12").Value = 161.75
$ws.Range("N12").Value = -872.71429

# Sheet CUL, row 104 (hunk 19)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3557
$ws.Range("J104").Value = 3557
$ws.Range("L104").Value = 10671
$ws.Range("N104").Value = -15913

# Sheet CUL, row 131 (hunk 20)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2701.9355
$ws.Range("I131").Value = 479
$ws.Range("J131").Value = 3129.423
$ws.Range("K131").Value = 1437
$ws.Range("L131").Value = 9388.269
$ws.Range("M131").Value = 3603
$ws.Range("N131").Value = -19468.269

# Sheet GSM, row 5 (hunk 21)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2144.5557
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 2966.8333
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 2966.8333
$ws.Range("M5").Value = -388
$ws.Range("N5").Value = -3190.8333

# Sheet GSM, row 80 (hunk 22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5996

# Sheet GSM, row 83 (hunk 23)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 20000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -29984

# Sheet GSM, row 132 (hunk 24)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2648.932
$ws.Range("I132").Value = 1883.6072
$ws.Range("J132").Value = 3988.25
$ws.Range("K132").Value = 5650.821599999999
$ws.Range("L132").Value = 11964.75
$ws.Range("M132").Value = -3120.821599999999
$ws.Range("N132").Value = -17024.75

# Sheet LTW, row 16 (hunk 25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2282
$ws.Range("I16").Value = 900
$ws.Range("J16").Value = 3664
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 3664
$ws.Range("M16").Value = -730
$ws.Range("N16").Value = -4004

# Sheet LTW, row 46 (hunk 26)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 144385.42
$ws.Range("I46").Value = 334300
$ws.Range("J46").Value = 1949.5
$ws.Range("K46").Value = 334300
$ws.Range("L46").Value = 1949.5
$ws.Range("M46").Value = -334112
$ws.Range("N46").Value = -2325.5

# Sheet LTW, row 61 (hunk 27)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2192.0833
$ws.Range("I61").Value = 2080
$ws.Range("J61").Value = 2272.1428
$ws.Range("K61").Value = 2080
$ws.Range("L61").Value = 2272.1428
$ws.Range("M61").Value = -1878
$ws.Range("N61").Value = -2676.1428

# Sheet LTW, row 93 (hunk 28)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 51001.5
$ws.Range("I93").Value = 51001.5
$ws.Range("K93").Value = 51001.5
$ws.Range("M93").Value = -49753.5

# Sheet LTW, row 100 (hunk 29)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2092.2856
$ws.Range("I100").Value = 1941
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1941
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1400
$ws.Range("N100").Value = -4082

# Sheet LTW, row 113 (hunk 30)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2192.0833
$ws.Range("I113").Value = 2080
$ws.Range("J113").Value = 2272.1428
$ws.Range("K113").Value = 2080
$ws.Range("L113").Value = 2272.1428
$ws.Range("M113").Value = 90
$ws.Range("N113").Value = -6612.1428
